# Update the "Förändrad" (changed) date column (C2:C54) from
# serial date 45183 (2023-09-14) to 45184 (2023-09-15) for every
# data row, matching the upstream automatic data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C54")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
